{"js": "// Update each lattice-multiplication exercise cell in the 5x3 table with\n// new problems, keeping the existing \"NN x NN / line2 / ---- / line4 / line5\"\n// five-line layout (separated by manual line breaks) and the run's 32\n// half-point font size.\nconst grid = [\n  [\n    [\"88 x 38\", \"  3    8\", \"  ----\", \"8|    |\", \"8|    |\"],\n    [\"74 x 73\", \"  7    3\", \"  ----\", \"7|    |\", \"4|    |\"],\n    [\"67 x 56\", \"  5    6\", \"  ----\", \"6|    |\", \"7|    |\"],\n  ],\n  [\n    [\"38 x 93\", \"  9    3\", \"  ----\", \"3|    |\", \"8|    |\"],\n    [\"55 x 42\", \"  4    2\", \"  ----\", \"5|    |\", \"5|    |\"],\n    [\"18 x 67\", \"  6    7\", \"  ----\", \"1|    |\", \"8|    |\"],\n  ],\n  [\n    [\"75 x 11\", \"  1    1\", \"  ----\", \"7|    |\", \"5|    |\"],\n    [\"29 x 34\", \"  3    4\", \"  ----\", \"2|    |\", \"9|    |\"],\n    [\"11 x 38\", \"  3    8\", \"  ----\", \"1|    |\", \"1|    |\"],\n  ],\n  [\n    [\"67 x 10\", \"  1    0\", \"  ----\", \"6|    |\", \"7|    |\"],\n    [\"74 x 26\", \"  2    6\", \"  ----\", \"7|    |\", \"4|    |\"],\n    [\"25 x 38\", \"  3    8\", \"  ----\", \"2|    |\", \"5|    |\"],\n  ],\n  [\n    [\"77 x 70\", \"  7    0\", \"  ----\", \"7|    |\", \"7|    |\"],\n    [\"88 x 43\", \"  4    3\", \"  ----\", \"8|    |\", \"8|    |\"],\n    [\"47 x 94\", \"  9    4\", \"  ----\", \"4|    |\", \"7|    |\"],\n  ],\n];\n\nfunction escapeXml(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\");\n}\n\n// Build a <w:p> fragment that mirrors the original cell content shape:\n// a single run at 32 half-points, with the 5 lines joined by <w:br/>,\n// keeping xml:space=\"preserve\" on any line that has leading/trailing\n// whitespace (exactly as Word itself would emit).\nfunction buildCellOoxml(lines) {\n  const runsXml = lines\n    .map((line, i) => {\n      const needsPreserve = /^\\s|\\s$/.test(line);\n      const preserveAttr = needsPreserve ? ' xml:space=\"preserve\"' : \"\";\n      const textXml = `<w:t${preserveAttr}>${escapeXml(line)}</w:t>`;\n      return i === 0 ? textXml : `<w:br/>${textXml}`;\n    })\n    .join(\"\");\n\n  return `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p><w:r><w:rPr><w:sz w:val=\"32\"/></w:rPr>${runsXml}</w:r></w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n}\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (let r = 0; r < grid.length; r++) {\n  for (let c = 0; c < grid[r].length; c++) {\n    const cell = table.getCell(r, c);\n    const ooxml = buildCellOoxml(grid[r][c]);\n    cell.body.insertOoxml(ooxml, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Update each lattice-multiplication exercise cell in the 5x3 table with\n# new problems, keeping the existing \"NN x NN / line2 / ---- / line4 / line5\"\n# five-line layout (separated by manual line breaks) and the run's 32\n# half-point font size.\n\nfunction Set-LatticeCell {\n    param(\n        $Table,\n        [int]$Row,\n        [int]$Col,\n        [string[]]$Lines\n    )\n\n    $cell = $Table.Cell($Row, $Col)\n\n    $runsXml = \"\"\n    for ($i = 0; $i -lt $Lines.Count; $i++) {\n        $line = $Lines[$i]\n        $needsPreserve = ($line -match '^\\s') -or ($line -match '\\s$')\n        $preserveAttr = \"\"\n        if ($needsPreserve) { $preserveAttr = ' xml:space=\"preserve\"' }\n        $escaped = $line -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'\n        $textXml = \"<w:t$preserveAttr>$escaped</w:t>\"\n        if ($i -eq 0) {\n            $runsXml += $textXml\n        } else {\n            $runsXml += \"<w:br/>$textXml\"\n        }\n    }\n\n    $frag = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body><w:p><w:r><w:rPr><w:sz w:val=\"32\"/></w:rPr>' + $runsXml + '</w:r></w:p></w:body>' +\n        '</w:document>' +\n        '</pkg:xmlData></pkg:part></pkg:package>'\n\n    [void]$cell.Range.InsertXML($frag)\n}\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\nSet-LatticeCell $table 1 1 @(\"88 x 38\", \"  3    8\", \"  ----\", \"8|    |\", \"8|    |\")\nSet-LatticeCell $table 1 2 @(\"74 x 73\", \"  7    3\", \"  ----\", \"7|    |\", \"4|    |\")\nSet-LatticeCell $table 1 3 @(\"67 x 56\", \"  5    6\", \"  ----\", \"6|    |\", \"7|    |\")\n\nSet-LatticeCell $table 2 1 @(\"38 x 93\", \"  9    3\", \"  ----\", \"3|    |\", \"8|    |\")\nSet-LatticeCell $table 2 2 @(\"55 x 42\", \"  4    2\", \"  ----\", \"5|    |\", \"5|    |\")\nSet-LatticeCell $table 2 3 @(\"18 x 67\", \"  6    7\", \"  ----\", \"1|    |\", \"8|    |\")\n\nSet-LatticeCell $table 3 1 @(\"75 x 11\", \"  1    1\", \"  ----\", \"7|    |\", \"5|    |\")\nSet-LatticeCell $table 3 2 @(\"29 x 34\", \"  3    4\", \"  ----\", \"2|    |\", \"9|    |\")\nSet-LatticeCell $table 3 3 @(\"11 x 38\", \"  3    8\", \"  ----\", \"1|    |\", \"1|    |\")\n\nSet-LatticeCell $table 4 1 @(\"67 x 10\", \"  1    0\", \"  ----\", \"6|    |\", \"7|    |\")\nSet-LatticeCell $table 4 2 @(\"74 x 26\", \"  2    6\", \"  ----\", \"7|    |\", \"4|    |\")\nSet-LatticeCell $table 4 3 @(\"25 x 38\", \"  3    8\", \"  ----\", \"2|    |\", \"5|    |\")\n\nSet-LatticeCell $table 5 1 @(\"77 x 70\", \"  7    0\", \"  ----\", \"7|    |\", \"7|    |\")\nSet-LatticeCell $table 5 2 @(\"88 x 43\", \"  4    3\", \"  ----\", \"8|    |\", \"8|    |\")\nSet-LatticeCell $table 5 3 @(\"47 x 94\", \"  9    4\", \"  ----\", \"4|    |\", \"7|    |\")\n"}
